$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "Mean"
$ws.Range("B5").Value = 0.001736295967875623
$ws.Range("C5").Value = 0.001875686121044729
$ws.Range("D5").Value = 0.02330513699821474
$ws.Range("E5").Value = 0.02382863760741438
$ws.Range("F5").Value = 0.009644119010622919
$ws.Range("G5").Value = 0.06538616969329449
